$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.252368
$ws.Cells.Item(2, 8).Value = 0.757104
$ws.Cells.Item(2, 9).Value = 0.1374452314236153
$ws.Cells.Item(2, 10).Value = 0.1374452314236153
$ws.Cells.Item(2, 13).Value = 0.2313123333333333
$ws.Cells.Item(2, 14).Value = 0.693937
$ws.Cells.Item(2, 15).Value = 0.7569517164947553
$ws.Cells.Item(2, 16).Value = 0.7569517164947555
$ws.Cells.Item(2, 17).Value = 0.05837583093866666
$ws.Cells.Item(2, 18).Value = 0.525382478448
$ws.Cells.Item(2, 19).Value = 0.1040394038501245
$ws.Cells.Item(2, 20).Value = 0.1040394038501245

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.252368
$ws.Cells.Item(3, 8).Value = 0.757104
$ws.Cells.Item(3, 9).Value = 0.1374452314236153
$ws.Cells.Item(3, 10).Value = 0.1374452314236153
$ws.Cells.Item(3, 15).Value = 0.2385552472206224
$ws.Cells.Item(3, 16).Value = 0.2385552472206224
$ws.Cells.Item(3, 17).Value = 0.01839729070933333
$ws.Cells.Item(3, 18).Value = 0.165575616384
$ws.Cells.Item(3, 19).Value = 0.03278828116155621
$ws.Cells.Item(3, 20).Value = 0.03278828116155622

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.252368
$ws.Cells.Item(4, 8).Value = 0.757104
$ws.Cells.Item(4, 9).Value = 0.1374452314236153
$ws.Cells.Item(4, 10).Value = 0.1374452314236153
$ws.Cells.Item(4, 13).Value = 0.001373
$ws.Cells.Item(4, 14).Value = 0.004119
$ws.Cells.Item(4, 15).Value = 0.004493036284622232
$ws.Cells.Item(4, 16).Value = 0.004493036284622232
$ws.Cells.Item(4, 17).Value = 0.000346501264
$ws.Cells.Item(4, 18).Value = 0.003118511376
$ws.Cells.Item(4, 19).Value = 0.0006175464119346035
$ws.Cells.Item(4, 20).Value = 0.0006175464119346035

# Row 5
$ws.Cells.Item(5, 9).Value = 0.1782163802407412
$ws.Cells.Item(5, 10).Value = 0.1782163802407412
$ws.Cells.Item(5, 13).Value = 0.2313123333333333
$ws.Cells.Item(5, 14).Value = 0.693937
$ws.Cells.Item(5, 15).Value = 0.7569517164947553
$ws.Cells.Item(5, 16).Value = 0.7569517164947555
$ws.Cells.Item(5, 17).Value = 0.07569218062844446
$ws.Cells.Item(5, 18).Value = 0.6812296256560001
$ws.Cells.Item(5, 19).Value = 0.1349011949307111
$ws.Cells.Item(5, 20).Value = 0.1349011949307111

# Row 6
$ws.Cells.Item(6, 9).Value = 0.1782163802407412
$ws.Cells.Item(6, 10).Value = 0.1782163802407412
$ws.Cells.Item(6, 15).Value = 0.2385552472206224
$ws.Cells.Item(6, 16).Value = 0.2385552472206224
$ws.Cells.Item(6, 19).Value = 0.04251445264709446
$ws.Cells.Item(6, 20).Value = 0.04251445264709446

# Row 7
$ws.Cells.Item(7, 9).Value = 0.1782163802407412
$ws.Cells.Item(7, 10).Value = 0.1782163802407412
$ws.Cells.Item(7, 13).Value = 0.001373
$ws.Cells.Item(7, 14).Value = 0.004119
$ws.Cells.Item(7, 15).Value = 0.004493036284622232
$ws.Cells.Item(7, 16).Value = 0.004493036284622232
$ws.Cells.Item(7, 17).Value = 0.0004492858746666667
$ws.Cells.Item(7, 18).Value = 0.004043572872
$ws.Cells.Item(7, 19).Value = 0.0008007326629356828
$ws.Cells.Item(7, 20).Value = 0.0008007326629356827

# Row 8
$ws.Cells.Item(8, 7).Value = 0.267684
$ws.Cells.Item(8, 8).Value = 0.803052
$ws.Cells.Item(8, 9).Value = 0.1457866660131199
$ws.Cells.Item(8, 10).Value = 0.1457866660131199
$ws.Cells.Item(8, 13).Value = 0.2313123333333333
$ws.Cells.Item(8, 14).Value = 0.693937
$ws.Cells.Item(8, 15).Value = 0.7569517164947553
$ws.Cells.Item(8, 16).Value = 0.7569517164947555
$ws.Cells.Item(8, 17).Value = 0.06191861063599999
$ws.Cells.Item(8, 18).Value = 0.557267495724
$ws.Cells.Item(8, 19).Value = 0.1103534670806787
$ws.Cells.Item(8, 20).Value = 0.1103534670806788

# Row 9
$ws.Cells.Item(9, 7).Value = 0.267684
$ws.Cells.Item(9, 8).Value = 0.803052
$ws.Cells.Item(9, 9).Value = 0.1457866660131199
$ws.Cells.Item(9, 10).Value = 0.1457866660131199
$ws.Cells.Item(9, 15).Value = 0.2385552472206224
$ws.Cells.Item(9, 16).Value = 0.2385552472206224
$ws.Cells.Item(9, 17).Value = 0.019513806688
$ws.Cells.Item(9, 18).Value = 0.175624260192
$ws.Cells.Item(9, 19).Value = 0.03477817415223013
$ws.Cells.Item(9, 20).Value = 0.03477817415223013

# Row 10
$ws.Cells.Item(10, 7).Value = 0.267684
$ws.Cells.Item(10, 8).Value = 0.803052
$ws.Cells.Item(10, 9).Value = 0.1457866660131199
$ws.Cells.Item(10, 10).Value = 0.1457866660131199
$ws.Cells.Item(10, 13).Value = 0.001373
$ws.Cells.Item(10, 14).Value = 0.004119
$ws.Cells.Item(10, 15).Value = 0.004493036284622232
$ws.Cells.Item(10, 16).Value = 0.004493036284622232
$ws.Cells.Item(10, 17).Value = 0.000367530132
$ws.Cells.Item(10, 18).Value = 0.003307771188
$ws.Cells.Item(10, 19).Value = 0.0006550247802110505
$ws.Cells.Item(10, 20).Value = 0.0006550247802110505

# Row 11
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.03422333333333333
$ws.Cells.Item(11, 8).Value = 0.10267
$ws.Cells.Item(11, 9).Value = 0.01863878926839984
$ws.Cells.Item(11, 10).Value = 0.01863878926839983
$ws.Cells.Item(11, 13).Value = 0.2313123333333333
$ws.Cells.Item(11, 14).Value = 0.693937
$ws.Cells.Item(11, 15).Value = 0.7569517164947553
$ws.Cells.Item(11, 16).Value = 0.7569517164947555
$ws.Cells.Item(11, 17).Value = 0.007916279087777778
$ws.Cells.Item(11, 18).Value = 0.07124651178999999
$ws.Cells.Item(11, 19).Value = 0.01410866353009928
$ws.Cells.Item(11, 20).Value = 0.01410866353009928

# Row 12
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.03422333333333333
$ws.Cells.Item(12, 8).Value = 0.10267
$ws.Cells.Item(12, 9).Value = 0.01863878926839984
$ws.Cells.Item(12, 10).Value = 0.01863878926839983
$ws.Cells.Item(12, 15).Value = 0.2385552472206224
$ws.Cells.Item(12, 16).Value = 0.2385552472206224
$ws.Cells.Item(12, 17).Value = 0.002494835368888889
$ws.Cells.Item(12, 18).Value = 0.02245351832
$ws.Cells.Item(12, 19).Value = 0.004446380981816205
$ws.Cells.Item(12, 20).Value = 0.004446380981816205

# Row 13
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.03422333333333333
$ws.Cells.Item(13, 8).Value = 0.10267
$ws.Cells.Item(13, 9).Value = 0.01863878926839984
$ws.Cells.Item(13, 10).Value = 0.01863878926839983
$ws.Cells.Item(13, 13).Value = 0.001373
$ws.Cells.Item(13, 14).Value = 0.004119
$ws.Cells.Item(13, 15).Value = 0.004493036284622232
$ws.Cells.Item(13, 16).Value = 0.004493036284622232
$ws.Cells.Item(13, 17).Value = 0.00004698863666666667
$ws.Cells.Item(13, 18).Value = 0.00042289773
$ws.Cells.Item(13, 19).Value = 0.00008374475648434793
$ws.Cells.Item(13, 20).Value = 0.0000837447564843479

# Row 14
$ws.Cells.Item(14, 7).Value = 0.6562846666666667
$ws.Cells.Item(14, 8).Value = 1.968854
$ws.Cells.Item(14, 9).Value = 0.3574272407348407
$ws.Cells.Item(14, 10).Value = 0.3574272407348407
$ws.Cells.Item(14, 13).Value = 0.2313123333333333
$ws.Cells.Item(14, 14).Value = 0.693937
$ws.Cells.Item(14, 15).Value = 0.7569517164947553
$ws.Cells.Item(14, 16).Value = 0.7569517164947555
$ws.Cells.Item(14, 17).Value = 0.1518067375775556
$ws.Cells.Item(14, 18).Value = 1.366260638198
$ws.Cells.Item(14, 19).Value = 0.2705551633962218
$ws.Cells.Item(14, 20).Value = 0.2705551633962218

# Row 15
$ws.Cells.Item(15, 7).Value = 0.6562846666666667
$ws.Cells.Item(15, 8).Value = 1.968854
$ws.Cells.Item(15, 9).Value = 0.3574272407348407
$ws.Cells.Item(15, 10).Value = 0.3574272407348407
$ws.Cells.Item(15, 15).Value = 0.2385552472206224
$ws.Cells.Item(15, 16).Value = 0.2385552472206224
$ws.Cells.Item(15, 17).Value = 0.04784227715377778
$ws.Cells.Item(15, 18).Value = 0.430580494384
$ws.Cells.Item(15, 19).Value = 0.08526614377688481
$ws.Cells.Item(15, 20).Value = 0.08526614377688482

# Row 16
$ws.Cells.Item(16, 7).Value = 0.6562846666666667
$ws.Cells.Item(16, 8).Value = 1.968854
$ws.Cells.Item(16, 9).Value = 0.3574272407348407
$ws.Cells.Item(16, 10).Value = 0.3574272407348407
$ws.Cells.Item(16, 13).Value = 0.001373
$ws.Cells.Item(16, 14).Value = 0.004119
$ws.Cells.Item(16, 15).Value = 0.004493036284622232
$ws.Cells.Item(16, 16).Value = 0.004493036284622232
$ws.Cells.Item(16, 17).Value = 0.0009010788473333335
$ws.Cells.Item(16, 18).Value = 0.008109709626000002
$ws.Cells.Item(16, 19).Value = 0.001605933561734045
$ws.Cells.Item(16, 20).Value = 0.001605933561734045

# Row 17
$ws.Cells.Item(17, 7).Value = 0.2983456666666667
$ws.Cells.Item(17, 8).Value = 0.8950370000000001
$ws.Cells.Item(17, 9).Value = 0.162485692319283
$ws.Cells.Item(17, 10).Value = 0.162485692319283
$ws.Cells.Item(17, 13).Value = 0.2313123333333333
$ws.Cells.Item(17, 14).Value = 0.693937
$ws.Cells.Item(17, 15).Value = 0.7569517164947553
$ws.Cells.Item(17, 16).Value = 0.7569517164947555
$ws.Cells.Item(17, 17).Value = 0.06901103229655556
$ws.Cells.Item(17, 18).Value = 0.621099290669
$ws.Cells.Item(17, 19).Value = 0.1229938237069199
$ws.Cells.Item(17, 20).Value = 0.12299382370692

# Row 18
$ws.Cells.Item(18, 7).Value = 0.2983456666666667
$ws.Cells.Item(18, 8).Value = 0.8950370000000001
$ws.Cells.Item(18, 9).Value = 0.162485692319283
$ws.Cells.Item(18, 10).Value = 0.162485692319283
$ws.Cells.Item(18, 15).Value = 0.2385552472206224
$ws.Cells.Item(18, 16).Value = 0.2385552472206224
$ws.Cells.Item(18, 17).Value = 0.02174900130577778
$ws.Cells.Item(18, 18).Value = 0.195741011752
$ws.Cells.Item(18, 19).Value = 0.03876181450104053
$ws.Cells.Item(18, 20).Value = 0.03876181450104053

# Row 19
$ws.Cells.Item(19, 7).Value = 0.2983456666666667
$ws.Cells.Item(19, 8).Value = 0.8950370000000001
$ws.Cells.Item(19, 9).Value = 0.162485692319283
$ws.Cells.Item(19, 10).Value = 0.162485692319283
$ws.Cells.Item(19, 13).Value = 0.001373
$ws.Cells.Item(19, 14).Value = 0.004119
$ws.Cells.Item(19, 15).Value = 0.004493036284622232
$ws.Cells.Item(19, 16).Value = 0.004493036284622232
$ws.Cells.Item(19, 17).Value = 0.0004096286003333334
$ws.Cells.Item(19, 18).Value = 0.003686657403
$ws.Cells.Item(19, 19).Value = 0.0007300541113225023
$ws.Cells.Item(19, 20).Value = 0.0007300541113225023

Write-Output "Updated Slit1-Robo2 LR-pair values with new TPM data"